# Update cryptocurrency price/volume data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.588.75'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.26%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.142.73'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.12%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '574.15'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.50%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '164.56'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.56%  '
$ws.Range("E7").Value = '  -0.11%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.576'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -5.28%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.157.13'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.90%  '
$ws.Range("E10").Value = '  -2.29%  '
$ws.Range("E11").Value = '  -2.35%  '
$ws.Range("E12").Value = '  -0.85%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.692.41'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.08%  '
$ws.Range("E14").Value = '  -1.75%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '64.580.89'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.16%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '25.08'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.68%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.153.77'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.18%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0000156'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.91%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '409.78'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.27%  '
$ws.Range("E20").Value = '  -2.28%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.50'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.69%  '
$ws.Range("E22").Value = '  -0.91%  '
$ws.Range("E23").Value = '  -0.10%  '
$ws.Range("E24").Value = '  -1.66%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.484'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.32%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.195'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.82%  '
$ws.Range("E27").Value = '  -2.60%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.90'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.04%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.995'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.16%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.08%  '
$ws.Range("E31").Value = '  -0.93%  '
$ws.Range("E32").Value = '  -2.23%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '163.47'
$ws.Range("D33").Style = "Normal"
$ws.Range("E34").Value = '  -3.22%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.30'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.72%  '
$ws.Range("E36").Value = '  +1.19%  '
$ws.Range("E37").Value = '  -0.11%  '
$ws.Range("E38").Value = '  -0.60%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.644.47'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.75%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '23.76'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.80%  '
$ws.Range("E41").Value = '  -2.92%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '38.22'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.64%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.691'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.47%  '
$ws.Range("E44").Value = '  -1.58%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.31'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.11%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '290.45'
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '21.39'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.10%  '
$ws.Range("E48").Value = '  -3.20%  '
$ws.Range("E49").Value = '  -0.18%  '
$ws.Range("E50").Value = '  -1.88%  '
$ws.Range("E51").Value = '  +0.51%  '
